# "Importar productos desde excel incluyendo el codigo producto"
# Add a new "Codigo" (product code) column as the new column A, shifting the
# existing Nombre/Descripcion/Precio/Proveedor/Serie/Tipo/Sucursal/Stock
# columns one place to the right (B..I), fill in the new product codes, and
# record who added the extra "Jujutsu Kaisen" serie entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A so existing columns (A..H) shift to (B..I)
$ws.Range("A1").EntireColumn.Insert()

# New "Codigo" column header and values
$ws.Range("A1").Value = "Codigo"
$ws.Range("A2").Value = "-"
$ws.Range("A3").Value = "KORIQ4zChClzm"

# Row 2's Serie (now column F after the shift) now also credits "juan"
$ws.Range("F2").Value = "Jujutsu Kaisen;juan"

$ws.Range("A4").Value = "KORIS5zZjCqzn"

# Match styling used by the rest of the table: header row style, data row style
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the new "Codigo" column a sensible custom width like the other columns
$ws.Range("A1").ColumnWidth = 20.6
